$d = $word.ActiveDocument

# Delete the 2nd and 3rd paragraphs (keep only the first paragraph)
$d.Paragraphs(3).Range.Delete()
$d.Paragraphs(2).Range.Delete()

# Replace the text of the first (remaining) paragraph
$d.Content.Find.Execute("Hello I’m here", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This file was modified", 2)
